$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 74, shifting existing rows 74-135 down to 75-136.
$ws.Rows("74:74").Insert()

# Populate the newly inserted row 74 with the new weekly price record.
$ws.Cells.Item(74,1).Value  = 9
$ws.Cells.Item(74,2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(74,3).Value  = "Metropolitana"
$ws.Cells.Item(74,4).Value  = 44818
$ws.Cells.Item(74,5).Value  = 13
$ws.Cells.Item(74,6).Value  = 100112022
$ws.Cells.Item(74,7).Value  = "Arveja Verde"
$ws.Cells.Item(74,8).Value  = "Perfection"
$ws.Cells.Item(74,9).Value  = "Primera"
$ws.Cells.Item(74,10).Value = 34
$ws.Cells.Item(74,11).Value = 27000
$ws.Cells.Item(74,12).Value = 29000
$ws.Cells.Item(74,13).Value = 28000
$ws.Cells.Item(74,14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(74,15).Value = "Provincia de Huasco"
$ws.Cells.Item(74,16).Value = 1120
$ws.Cells.Item(74,17).Value = 25
$ws.Cells.Item(74,18).Value = "Hortaliza"
